# Add image scale (pixel limit) and colormap (contrast/bias) parameter
# columns J:Y to the dataset sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("J1").Value = "I4 Low Pixel Limit"
$ws.Range("K1").Value = "I4 High Pixel Limit"
$ws.Range("L1").Value = "I2 Low Pixel Limit"
$ws.Range("M1").Value = "I2 High Pixel Limit"
$ws.Range("N1").Value = "I1 Low Pixel Limit"
$ws.Range("O1").Value = "I1 High Pixel Limit"
$ws.Range("P1").Value = "I3 Low Pixel Limit"
$ws.Range("Q1").Value = "I3 High Pixel Limit"
$ws.Range("R1").Value = "I4 Contrast"
$ws.Range("S1").Value = "I4 Bias"
$ws.Range("T1").Value = "I2 Contrast"
$ws.Range("U1").Value = "I2 Bias"
$ws.Range("V1").Value = "I1 Contrast"
$ws.Range("W1").Value = "I1 Bias"
$ws.Range("X1").Value = "I3 Contrast"
$ws.Range("Y1").Value = "I3 Bias"

# Data row (row 2)
$ws.Range("J2").Value = 5.63551
$ws.Range("K2").Value = 11.0478
$ws.Range("L2").Value = -0.170264
$ws.Range("M2").Value = 4.11561
$ws.Range("N2").Value = -0.320132
$ws.Range("O2").Value = 3.357
$ws.Range("P2").Value = -0.320132
$ws.Range("Q2").Value = 14.1618
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 0.5
$ws.Range("T2").Value = 0.993151
$ws.Range("U2").Value = 0.4520547945205479
$ws.Range("V2").Value = 1.09589
$ws.Range("W2").Value = 0.5
$ws.Range("X2").Value = 1.23288
$ws.Range("Y2").Value = 0.414384
